# "fusion esteban / principal"
# Append the new "ASTRID MONIQUE" record as row 16 of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "ASTRID"
$ws.Range("B16").Value = "MONIQUE"
$ws.Range("C16").Value = "30 Rue Rambaud, 17000 La Rochelle"
$ws.Range("D16").Value = 46.16360619065998
$ws.Range("E16").Value = -1.155014376554859
